$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new blank row at row 43, shifting existing rows 43.. down by one.
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row 43 with the new transaction data.
$ws.Range("R43").Value = "value discovery debit icici"
$ws.Range("S43").Value = "2024-09-19 14:34:40"
